# Applies the Jenova_Profits market-data refresh across all 8 job sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR). Each block updates the
# currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ /
# LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ columns (H:N)
# for a specific Leve row, as refreshed by the scheduled market-data
# runner. A few rows gain a newly-computed LeveProfitNQ (M) cell, and two
# rows in CUL lose their LeveProfitHQ (N) cell (HQ price no longer tracked).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(18, 8).Value = 573.3333
$ws.Cells.Item(18, 9).Value = 433.25
$ws.Cells.Item(18, 10).Value = 1133.6666
$ws.Cells.Item(18, 11).Value = 433.25
$ws.Cells.Item(18, 12).Value = 1133.6666
$ws.Cells.Item(18, 13).Value = -149.25
$ws.Cells.Item(18, 14).Value = -1701.6666

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(32, 8).Value = 6039.5
$ws.Cells.Item(32, 9).Value = 4332.3335
$ws.Cells.Item(32, 10).Value = 6771.143
$ws.Cells.Item(32, 11).Value = 4332.3335
$ws.Cells.Item(32, 12).Value = 6771.143
$ws.Cells.Item(32, 13).Value = -4006.3335
$ws.Cells.Item(32, 14).Value = -7423.143

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(62, 8).Value = 4466347
$ws.Cells.Item(62, 9).Value = 5953929
$ws.Cells.Item(62, 10).Value = 3600.2856
$ws.Cells.Item(62, 11).Value = 5953929
$ws.Cells.Item(62, 12).Value = 3600.2856
$ws.Cells.Item(62, 13).Value = -5953305
$ws.Cells.Item(62, 14).Value = -4848.2856

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(65, 8).Value = 4466347
$ws.Cells.Item(65, 9).Value = 5953929
$ws.Cells.Item(65, 10).Value = 3600.2856
$ws.Cells.Item(65, 11).Value = 29769645
$ws.Cells.Item(65, 12).Value = 18001.428
$ws.Cells.Item(65, 13).Value = -29766525
$ws.Cells.Item(65, 14).Value = -24241.428

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(116, 8).Value = 4971.6787
$ws.Cells.Item(116, 9).Value = 4801.4287
$ws.Cells.Item(116, 11).Value = 4801.4287
$ws.Cells.Item(116, 13).Value = -1359.4287

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(132, 8).Value = 2239.3572
$ws.Cells.Item(132, 9).Value = 1868.8695
$ws.Cells.Item(132, 11).Value = 5606.6085
$ws.Cells.Item(132, 13).Value = -3076.6085

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 123945.445
$ws.Cells.Item(2, 9).Value = 158644.28
$ws.Cells.Item(2, 11).Value = 158644.28
$ws.Cells.Item(2, 13).Value = -158531.28

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(4, 8).Value = 9
$ws.Cells.Item(4, 9).Value = 8
$ws.Cells.Item(4, 11).Value = 8
$ws.Cells.Item(4, 13).Value = 108

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 3769.6853
$ws.Cells.Item(32, 9).Value = 2629.0852
$ws.Cells.Item(32, 11).Value = 2629.0852
$ws.Cells.Item(32, 13).Value = -2342.0852

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 4244.5835
$ws.Cells.Item(61, 9).Value = 2816.9524
$ws.Cells.Item(61, 11).Value = 2816.9524
$ws.Cells.Item(61, 13).Value = -2604.9524

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(116, 8).Value = 123945.445
$ws.Cells.Item(116, 9).Value = 158644.28
$ws.Cells.Item(116, 11).Value = 158644.28
$ws.Cells.Item(116, 13).Value = -156350.28

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(136, 8).Value = 4244.5835
$ws.Cells.Item(136, 9).Value = 2816.9524
$ws.Cells.Item(136, 11).Value = 8450.8572
$ws.Cells.Item(136, 13).Value = -5900.8572

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 123945.445
$ws.Cells.Item(3, 9).Value = 158644.28
$ws.Cells.Item(3, 11).Value = 158644.28
$ws.Cells.Item(3, 13).Value = -158530.28

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(7, 8).Value = 1400
$ws.Cells.Item(7, 9).Value = 500
$ws.Cells.Item(7, 10).Value = 1625
$ws.Cells.Item(7, 11).Value = 500
$ws.Cells.Item(7, 12).Value = 1625
$ws.Cells.Item(7, 13).Value = -387
$ws.Cells.Item(7, 14).Value = -1851

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 2376.9355
$ws.Cells.Item(94, 9).Value = 764.2174
$ws.Cells.Item(94, 10).Value = 7013.5
$ws.Cells.Item(94, 11).Value = 764.2174
$ws.Cells.Item(94, 12).Value = 7013.5
$ws.Cells.Item(94, 13).Value = -313.2174
$ws.Cells.Item(94, 14).Value = -7915.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 16647.328
$ws.Cells.Item(134, 9).Value = 1910.8422
$ws.Cells.Item(134, 11).Value = 5732.5266
$ws.Cells.Item(134, 13).Value = -3197.5266

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 203760.12
$ws.Cells.Item(31, 9).Value = 272649.72
$ws.Cells.Item(31, 11).Value = 272649.72
$ws.Cells.Item(31, 13).Value = -272354.72

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(32, 8).Value = 1580
$ws.Cells.Item(32, 9).Value = 725
$ws.Cells.Item(32, 10).Value = 5000
$ws.Cells.Item(32, 11).Value = 725
$ws.Cells.Item(32, 12).Value = 5000
$ws.Cells.Item(32, 13).Value = -409
$ws.Cells.Item(32, 14).Value = -5632

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(33, 8).Value = 1500
$ws.Cells.Item(33, 9).Value = 1500
$ws.Cells.Item(33, 11).Value = 1500
$ws.Cells.Item(33, 13).Value = -1121

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(34, 8).Value = 203760.12
$ws.Cells.Item(34, 9).Value = 272649.72
$ws.Cells.Item(34, 11).Value = 272649.72
$ws.Cells.Item(34, 13).Value = -272447.72

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(35, 8).Value = 878.9
$ws.Cells.Item(35, 9).Value = 973.44446
$ws.Cells.Item(35, 11).Value = 973.44446
$ws.Cells.Item(35, 13).Value = -679.44446

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(36, 8).Value = 5047.5
$ws.Cells.Item(36, 9).Value = 5047.5
$ws.Cells.Item(36, 11).Value = 5047.5
$ws.Cells.Item(36, 13).Value = -4659.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(40, 8).Value = 5047.5
$ws.Cells.Item(40, 9).Value = 5047.5
$ws.Cells.Item(40, 11).Value = 5047.5
$ws.Cells.Item(40, 13).Value = -4887.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(105, 8).Value = 2098.125
$ws.Cells.Item(105, 9).Value = 1943
$ws.Cells.Item(105, 11).Value = 1943
$ws.Cells.Item(105, 13).Value = -196

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(134, 8).Value = 203619.72
$ws.Cells.Item(134, 9).Value = 2411.4688
$ws.Cells.Item(134, 11).Value = 7234.4064
$ws.Cells.Item(134, 13).Value = -4699.4064

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(23, 8).Value = 277.27274
$ws.Cells.Item(23, 10).Value = 116.666664
$ws.Cells.Item(23, 12).Value = 349.999992
$ws.Cells.Item(23, 14).Value = -819.999992

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(64, 8).Value = 500000350
$ws.Cells.Item(64, 9).Value = 500000350
$ws.Cells.Item(64, 10).Value = 0
$ws.Cells.Item(64, 11).Value = 1500001050
$ws.Cells.Item(64, 12).Value = 0
$ws.Cells.Item(64, 14).ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(67, 8).Value = 500000350
$ws.Cells.Item(67, 9).Value = 500000350
$ws.Cells.Item(67, 10).Value = 0
$ws.Cells.Item(67, 11).Value = 1500001050
$ws.Cells.Item(67, 12).Value = 0
$ws.Cells.Item(67, 14).ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 1931.7368
$ws.Cells.Item(97, 9).Value = 1891.0769
$ws.Cells.Item(97, 11).Value = 1891.0769
$ws.Cells.Item(97, 13).Value = -1395.0769

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(113, 8).Value = 1252163.2
$ws.Cells.Item(113, 9).Value = 1667884.4
$ws.Cells.Item(113, 11).Value = 1667884.4
$ws.Cells.Item(113, 13).Value = -1665714.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 3441.2307
$ws.Cells.Item(122, 9).Value = 2023.7
$ws.Cells.Item(122, 10).Value = 8166.3335
$ws.Cells.Item(122, 11).Value = 6071.1
$ws.Cells.Item(122, 12).Value = 24499.0005
$ws.Cells.Item(122, 13).Value = -3621.1
$ws.Cells.Item(122, 14).Value = -29399.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 1214.2307
$ws.Cells.Item(16, 9).Value = 1214.2307
$ws.Cells.Item(16, 11).Value = 1214.2307
$ws.Cells.Item(16, 13).Value = -1044.2307

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(32, 8).Value = 5125
$ws.Cells.Item(32, 9).Value = 5125
$ws.Cells.Item(32, 11).Value = 5125
$ws.Cells.Item(32, 13).Value = -4808

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(68, 8).Value = 94950.09
$ws.Cells.Item(68, 9).Value = 3466.8333
$ws.Cells.Item(68, 11).Value = 3466.8333
$ws.Cells.Item(68, 13).Value = -2717.8333

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(71, 8).Value = 94950.09
$ws.Cells.Item(71, 9).Value = 3466.8333
$ws.Cells.Item(71, 11).Value = 17334.1665
$ws.Cells.Item(71, 13).Value = -13590.1665

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(100, 8).Value = 1826.8572
$ws.Cells.Item(100, 9).Value = 1577.8
$ws.Cells.Item(100, 11).Value = 1577.8
$ws.Cells.Item(100, 13).Value = -1036.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(122, 8).Value = 609965.3
$ws.Cells.Item(122, 9).Value = 438128.75
$ws.Cells.Item(122, 11).Value = 1314386.25
$ws.Cells.Item(122, 13).Value = -1311936.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(2, 8).Value = 44428.57
$ws.Cells.Item(2, 9).Value = 44428.57
$ws.Cells.Item(2, 11).Value = 44428.57
$ws.Cells.Item(2, 13).Value = -44316.57

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(33, 8).Value = 5006.3335
$ws.Cells.Item(33, 9).Value = 19
$ws.Cells.Item(33, 11).Value = 19
$ws.Cells.Item(33, 13).Value = 231

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(36, 8).Value = 5006.3335
$ws.Cells.Item(36, 9).Value = 19
$ws.Cells.Item(36, 11).Value = 19
$ws.Cells.Item(36, 13).Value = 231

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 110999.8
$ws.Cells.Item(62, 9).Value = 209999.6
$ws.Cells.Item(62, 11).Value = 209999.6
$ws.Cells.Item(62, 13).Value = -209375.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(65, 8).Value = 110999.8
$ws.Cells.Item(65, 9).Value = 209999.6
$ws.Cells.Item(65, 11).Value = 1049998
$ws.Cells.Item(65, 13).Value = -1046878

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(96, 8).Value = 115144.78
$ws.Cells.Item(96, 9).Value = 203440
$ws.Cells.Item(96, 10).Value = 4775.75
$ws.Cells.Item(96, 11).Value = 203440
$ws.Cells.Item(96, 12).Value = 4775.75
$ws.Cells.Item(96, 13).Value = -202067
$ws.Cells.Item(96, 14).Value = -7521.75
